$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: shorten the "is not a valid key..." sentence and replace its tail
# with "... and as such this fails BCNF."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " is not a valid key, for it to be valid Engine_Size would have to be a superkey in CAR_DETAILS1.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " is not a valid key and as such this fails BCNF.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 2: CAR_DETAILS1's attribute list re-ordered from
# "Registration_No, Make, Model, Colour)" to "Registration_No, Colour, Make, Model)"
# The underlined "Registration_No" run is left untouched; only the plain-text
# tail after it is rewritten.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    ", Make, Model, Colour)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", Colour, Make, Model)", 2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 3: CAR_DETAILS3's key changes from {Registration_No, Engine_Size} to
# {Make, Model}, with Engine_Size staying in the relation but losing its
# underline (key) formatting, and Make/Model gaining it.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("CAR_DETAILS3 (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos = $anchor.End

# "Registration_No" (underlined) -> "Make" (keep underline)
$rMake = $d.Range($pos, $pos + 15)
$rMake.Text = "Make"
$pos = $rMake.End

# "*, " (plain) -> ", " (plain) -- drop the old key-marking asterisk
$rComma = $d.Range($pos, $pos + 3)
$rComma.Text = ", "
$pos = $rComma.End

# insert new underlined "Model"
$rModel = $d.Range($pos, $pos)
$rModel.InsertAfter("Model")
$rModel.Font.Underline = 1
$pos = $rModel.End

# insert the new key-marking asterisk + separator (plain)
$rStar = $d.Range($pos, $pos)
$rStar.InsertAfter("*, ")
$rStar.Font.Underline = 0
$pos = $rStar.End

# "Engine_Size" loses its underline (no longer part of the key)
$rEngine = $d.Range($pos, $pos + 11)
$rEngine.Font.Underline = 0

# ---------------------------------------------------------------------------
# Edit 4: split the trailing run of blank lines into two paragraphs -- an
# extra empty paragraph is inserted after the second trailing line break.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$paraText = $lastPara.Range.Text
$splitOffset = $paraText.IndexOf([char]11, [char]11, 0)

# locate the end of the 2nd manual line break (vertical-tab, chr(11)) run
$vtCount = 0
$relPos = -1
for ($i = 0; $i -lt $paraText.Length; $i++) {
    if ([int][char]$paraText[$i] -eq 11) {
        $vtCount++
        if ($vtCount -eq 2) {
            $relPos = $i + 1
            break
        }
    }
}

$splitPos = $lastPara.Range.Start + $relPos
$r = $d.Range($splitPos, $splitPos)
$r.InsertParagraphAfter()

$newPara = $d.Paragraphs($lastPara.Index + 1)
$r2 = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$r2.InsertParagraphBefore()
